$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (column D) and 1h volume/change (column E) values.
# Column D values are textual (not numeric) in the source data, even when they look like
# numbers (e.g. '313.82'), so a leading apostrophe forces Excel to store them as text,
# matching the original inline-string cell type instead of auto-coercing to a Number.

$ws.Range("D2").Value = "'27.379.03"
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = "'1.779.39"
$ws.Range("E3").Value = '  +3.73%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'313.82"
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = "'0.5316"
$ws.Range("E7").Value = '  +11.96%  '
$ws.Range("D8").Value = "'0.3762"
$ws.Range("E8").Value = '  +9.07%  '
$ws.Range("D9").Value = "'42.83"
$ws.Range("E9").Value = '  +1.49%  '
$ws.Range("D10").Value = "'0.07409"
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").Value = "'20.69"
$ws.Range("E13").Value = '  +4.65%  '
$ws.Range("D14").Value = "'6.108"
$ws.Range("E14").Value = '  +4.92%  '
$ws.Range("D15").Value = "'1.783.64"
$ws.Range("E15").Value = '  +3.76%  '
$ws.Range("E16").Value = '  +2.62%  '
$ws.Range("D17").Value = "'89.76"
$ws.Range("E17").Value = '  +2.98%  '
$ws.Range("D18").Value = "'0.00001056"
$ws.Range("E18").Value = '  +2.07%  '
$ws.Range("D19").Value = "'0.06434"
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = "'16.79"
$ws.Range("E21").Value = '  +2.19%  '
$ws.Range("D22").Value = "'5.907"
$ws.Range("D23").Value = "'27.428.97"
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").Value = "'2.096"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'155.55"
$ws.Range("E26").Value = '  +3.00%  '
$ws.Range("D27").Value = "'20.22"
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").Value = "'2.368"
$ws.Range("E28").Value = '  +15.01%  '
$ws.Range("D29").Value = "'1.988.32"
$ws.Range("E29").Value = '  +3.92%  '
$ws.Range("D30").Value = "'121.40"
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").Value = "'1.086"
$ws.Range("E31").Value = '  +6.11%  '
$ws.Range("D32").Value = "'0.1028"
$ws.Range("E32").Value = '  +11.77%  '
$ws.Range("D33").Value = "'5.589"
$ws.Range("E33").Value = '  +5.22%  '
$ws.Range("D34").Value = "'3.634"
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("D35").Value = "'0.02257"
$ws.Range("E35").Value = '  +3.65%  '
$ws.Range("D36").Value = "'0.05977"
$ws.Range("D37").Value = "'11.29"
$ws.Range("E37").Value = '  +3.21%  '
$ws.Range("D38").Value = "'4.918"
$ws.Range("E38").Value = '  +4.46%  '
$ws.Range("D39").Value = "'0.2052"
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("D40").Value = "'0.6132"
$ws.Range("E40").Value = '  +2.90%  '
$ws.Range("E41").Value = '  +10.06%  '
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("D43").Value = "'1.133"
$ws.Range("E43").Value = '  +4.57%  '
$ws.Range("D44").Value = "'13.22"
$ws.Range("E44").Value = '  +3.91%  '
$ws.Range("D45").Value = "'0.5794"
$ws.Range("E45").Value = '  +4.11%  '
$ws.Range("D46").Value = "'3.627"
$ws.Range("E46").Value = '  +1.30%  '
$ws.Range("D47").Value = "'121.52"
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("D48").Value = "'1.896"
$ws.Range("E48").Value = '  +3.85%  '
$ws.Range("D49").Value = "'1.120"
$ws.Range("E49").Value = '  +1.30%  '
$ws.Range("D50").Value = "'0.06729"
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("D51").Value = "'70.97"
$ws.Range("E51").Value = '  +2.63%  '
